$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("2022-01-11", "overview", "K02000001", "United Kingdom", 14732594, 120821, 379, 150609),
    @("2022-01-12", "overview", "K02000001", "United Kingdom", 14862138, 129587, 398, 151007),
    @("2022-01-13", "overview", "K02000001", "United Kingdom", 14967817, 109133, 335, 151342),
    @("2022-01-14", "overview", "K02000001", "United Kingdom", 15066395, 99652, 270, 151612),
    @("2022-01-15", "overview", "K02000001", "United Kingdom", 15147120, 81713, 287, 151899),
    @("2022-01-16", "overview", "K02000001", "United Kingdom", 15217280, 70924, 88, 151987)
)

$startRow = 518
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $row[0]
    $dateCell.ClearFormats()

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}
